# "Updated: - risk influence - took out supply chain"
# => take this version for interim presentation
#
# This script:
#   1. Updates a handful of numeric estimates (supply-chain investment /
#      hay-harvest labor rows) that were re-evaluated.
#   2. Strips the redundant applyFill/applyBorder cell formats that Excel
#      had minted for a few cells in the "risk influence" table, reverting
#      them back to the plain border-only formats already used by their
#      neighboring cells (same visual border, just without the spurious
#      applyFill/applyBorder flags).
#   3. Moves the sheet's scroll/selection down to the risk-influence block
#      that is now the focus of the interim presentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Numeric estimate updates -----------------------------------------

$ws.Range("C34").Value = 0.25
$ws.Range("E34").Value = 0.5

$ws.Range("C38").Value = 0.25
$ws.Range("E38").Value = 0.5

$ws.Range("E39").Value = 5

$ws.Range("E41").Value = 10000

# --- 2. Drop the redundant applyFill/applyBorder cell formats -------------
# (xlEdgeLeft = 7, xlEdgeRight = 10, xlLineStyleNone = -4142)

# Cells that only had a border "apply" flag but no actual border -> plain
# default format (matches the sibling cells around them).
$noBorderCells = @("G24", "G44")
foreach ($addr in $noBorderCells) {
    $ws.Range($addr).Borders.LineStyle = -4142
}

# Cells whose real border is "thin right edge only" (same border already
# used elsewhere as the plain, non-applyFill/applyBorder format).
$rightBorderCells = @("A50", "E50", "E51", "A54", "E54", "A55", "E55", "A58", "E58", "A59", "E59")
foreach ($addr in $rightBorderCells) {
    $ws.Range($addr).Borders.Item(10).LineStyle = 1
}

# Cells whose real border is "thin left + right edges" (same border already
# used elsewhere as the plain, non-applyFill/applyBorder format).
$leftRightBorderCells = @("F50", "F51", "F53", "F54", "F57", "F58")
foreach ($addr in $leftRightBorderCells) {
    $ws.Range($addr).Borders.Item(7).LineStyle = 1
    $ws.Range($addr).Borders.Item(10).LineStyle = 1
}

# --- 3. Move the view/selection to the risk-influence block ---------------

$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E42").Select() | Out-Null
